$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 25
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "06/22/2025"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 0.0004844799999999989
$ws.Cells.Item($row, 3).Value = 103203.4346103041
$ws.Cells.Item($row, 4).Value = 50
